$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "Lori"
$ws.Range("B20").Value = "Lightfoot"
$ws.Range("A20:B20").Font.Color = $ws.Range("A19:B19").Font.Color
